# This workbook has three sheets (Q4_18_19, Q4_17_18, Q4_16_17), each a small
# Group/Project/Total Forecast table with 5 data rows (rows 2-6). The edit
# trims every sheet down to a single data row (row 2) with updated values,
# and changes which row's cell-fill formatting survives.

$wb = $excel.ActiveWorkbook

$wsCurrent  = $wb.Worksheets.Item(1)   # Q4_18_19
$wsPrior    = $wb.Worksheets.Item(2)   # Q4_17_18
$wsTwoPrior = $wb.Worksheets.Item(3)   # Q4_16_17

# --- Re-point the formatting before we touch values/rows -------------------
# Q4_16_17!C2 currently carries the red highlight fill; Q4_18_19!C2 and
# Q4_17_18!C2 need that same red fill afterwards, and Q4_16_17!C2 needs to
# lose it. Copy/PasteSpecial(formats) reuses the workbook's existing style
# slot instead of minting a fresh (duplicate) one, so the orange fill that
# only Q4_17_18!C2 used stops being referenced by anything.
$wsTwoPrior.Range("C2").Copy()
$wsCurrent.Range("C2").PasteSpecial(-4122)
$wsPrior.Range("C2").PasteSpecial(-4122)
$wsTwoPrior.Range("C2").Style = "Normal"

# --- Q4_18_19: keep only the HSMRPG/A11 row ---------------------------------
$wsCurrent.Rows("3:6").Delete()
$wsCurrent.Range("A2").Value = "HSMRPG"
$wsCurrent.Range("B2").Value = "A11"
$wsCurrent.Range("C2").Value = 2739.7

# --- Q4_17_18: keep only the HSMRPG/A11 row ---------------------------------
$wsPrior.Rows("3:6").Delete()
$wsPrior.Range("A2").Value = "HSMRPG"
$wsPrior.Range("B2").Value = "A11"
$wsPrior.Range("C2").Value = 902

# --- Q4_16_17: keep only the Rail Group/A11 row -----------------------------
$wsTwoPrior.Rows("3:6").Delete()
$wsTwoPrior.Range("A2").Value = "Rail Group"
$wsTwoPrior.Range("B2").Value = "A11"
$wsTwoPrior.Range("C2").Value = 1378.6
